# Append 6 more players (18 rows: Group1/Group2/Difference each) to the
# LB aggregate sheet, continuing the existing alternating row-style pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Eric Wilson","Group1",2.888888888888889,71.55555555555556,37.11111111111111,34.44444444444444),
    @("Eric Wilson","Group2",0.6666666666666666,40,21,19),
    @("Eric Wilson","Difference",-2.222222222222222,-31.55555555555556,-16.11111111111111,-15.44444444444444),
    @("Jadeveon Clowney","Group1",3,29,19.66666666666667,9.333333333333334),
    @("Jadeveon Clowney","Group2",4,39,20.66666666666667,18.33333333333333),
    @("Jadeveon Clowney","Difference",1,10,1,8.999999999999998),
    @("Jihad Ward","Group1",0.8888888888888888,17.55555555555556,6.333333333333333,11.22222222222222),
    @("Jihad Ward","Group2",1.666666666666667,25.66666666666667,16.33333333333333,9.333333333333334),
    @("Jihad Ward","Difference",0.7777777777777779,8.111111111111111,10,-1.888888888888888),
    @("Kamu Grugier-Hill","Group1",1,53,36.33333333333334,16.66666666666667),
    @("Kamu Grugier-Hill","Group2",1.888888888888889,35.33333333333334,22.44444444444444,12.88888888888889),
    @("Kamu Grugier-Hill","Difference",0.8888888888888886,-17.66666666666666,-13.88888888888889,-3.777777777777777),
    @("Mack Wilson","Group1",3,54.33333333333334,35.66666666666666,18.66666666666667),
    @("Mack Wilson","Group2",3,49.33333333333334,27.66666666666667,21.66666666666667),
    @("Mack Wilson","Difference",0,-5,-7.999999999999996,3),
    @("Oren Burks","Group1",0,22.66666666666667,15,7.666666666666667),
    @("Oren Burks","Group2",1,41.66666666666666,21.33333333333333,20.33333333333333),
    @("Oren Burks","Difference",1,19,6.333333333333332,12.66666666666666)
)

$startRow = 14
$endRow = $startRow + $data.Count - 1

# Pre-stamp formatting for the new block by copying the existing alternating
# 3-row (Group1/Group2/Difference) style blocks (rows 2-4 = style "s=2",
# rows 5-7 = style "s=3") down across the new rows, continuing the pattern.
$srcGreen = $ws.Range("A2:F4")
$srcYellow = $ws.Range("A5:F7")
for ($r = $startRow; $r -le $endRow; $r += 6) {
    $srcGreen.Copy()
    $ws.Range("A" + $r + ":F" + ($r + 2)).PasteSpecial(-4122)
    $srcYellow.Copy()
    $ws.Range("A" + ($r + 3) + ":F" + ($r + 5)).PasteSpecial(-4122)
}

# Now write the values into the newly formatted rows.
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
}

Write-Output "done"
